$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 104 date label as text (shared string), avoiding Excel's
# automatic date auto-conversion of date-like strings by temporarily
# formatting the cell as Text, then reverting the number format so the
# cell keeps no explicit style (matching the rest of column A).
$dateCell = $ws.Range("A104")
$dateCell.NumberFormat = "@"
$dateCell.Value = "01-07-2021"
$dateCell.Style = "Normal"

# Updated / new numeric values for quarters affected by the reseasonalized
# Imacec series refresh (rows 50-104, columns B:J).
$ws.Range("H50").Value = 79.40000000000001
$ws.Range("I50").Value = 83.7
$ws.Range("D52").Value = 90.8
$ws.Range("H53").Value = 82.09999999999999
$ws.Range("I53").Value = 82.3
$ws.Range("B54").Value = 80.90000000000001
$ws.Range("D54").Value = 89.5
$ws.Range("D55").Value = 92.2
$ws.Range("H55").Value = 80.8
$ws.Range("C56").Value = 87.40000000000001
$ws.Range("D56").Value = 94
$ws.Range("G56").Value = 64.40000000000001
$ws.Range("H57").Value = 82.3
$ws.Range("H58").Value = 83.5
$ws.Range("J58").Value = 80.8
$ws.Range("D60").Value = 97.90000000000001
$ws.Range("J61").Value = 87.2
$ws.Range("D63").Value = 91.40000000000001
$ws.Range("D64").Value = 87.90000000000001
$ws.Range("H65").Value = 93.5
$ws.Range("D66").Value = 93.7
$ws.Range("D68").Value = 93.5
$ws.Range("H68").Value = 97.40000000000001
$ws.Range("F69").Value = 96.59999999999999
$ws.Range("D72").Value = 102.7
$ws.Range("F72").Value = 100.3
$ws.Range("F73").Value = 100.6
$ws.Range("D74").Value = 103
$ws.Range("E75").Value = 100.3
$ws.Range("H75").Value = 102
$ws.Range("D76").Value = 101.4
$ws.Range("D77").Value = 101.5
$ws.Range("F79").Value = 104.9
$ws.Range("D80").Value = 96.3
$ws.Range("H80").Value = 106.3
$ws.Range("J80").Value = 104.9
$ws.Range("D81").Value = 100.5
$ws.Range("F81").Value = 105.9
$ws.Range("I81").Value = 104.8
$ws.Range("D82").Value = 103.4
$ws.Range("F82").Value = 106
$ws.Range("H82").Value = 107.8
$ws.Range("B83").Value = 105.5
$ws.Range("D83").Value = 100.1
$ws.Range("F83").Value = 107.3
$ws.Range("D84").Value = 96.59999999999999
$ws.Range("H84").Value = 109
$ws.Range("J84").Value = 106.8
$ws.Range("D85").Value = 96.09999999999999
$ws.Range("F85").Value = 107.5
$ws.Range("B86").Value = 105.1
$ws.Range("F86").Value = 107
$ws.Range("H86").Value = 108.3
$ws.Range("J86").Value = 106.7
$ws.Range("D87").Value = 96.7
$ws.Range("F87").Value = 105.6
$ws.Range("H87").Value = 109.3
$ws.Range("I87").Value = 106.2
$ws.Range("B88").Value = 108.2
$ws.Range("D88").Value = 102
$ws.Range("F88").Value = 108
$ws.Range("H88").Value = 110.9
$ws.Range("J88").Value = 108.9
$ws.Range("B89").Value = 109
$ws.Range("D89").Value = 103.4
$ws.Range("F89").Value = 108.4
$ws.Range("G89").Value = 110.7
$ws.Range("H89").Value = 111.1
$ws.Range("J89").Value = 109.6
$ws.Range("B90").Value = 110.3
$ws.Range("C90").Value = 106.1
$ws.Range("F90").Value = 110.1
$ws.Range("H90").Value = 112.4
$ws.Range("J90").Value = 110.7
$ws.Range("B91").Value = 111.6
$ws.Range("C91").Value = 106.3
$ws.Range("D91").Value = 101.3
$ws.Range("F91").Value = 110.9
$ws.Range("H91").Value = 114.1
$ws.Range("I91").Value = 111.2
$ws.Range("J91").Value = 112.6
$ws.Range("C92").Value = 105
$ws.Range("D92").Value = 100.9
$ws.Range("F92").Value = 109.1
$ws.Range("H92").Value = 114.2
$ws.Range("I92").Value = 110.7
$ws.Range("J92").Value = 112.1
$ws.Range("B93").Value = 112
$ws.Range("C93").Value = 106.4
$ws.Range("D93").Value = 103.6
$ws.Range("F93").Value = 110.2
$ws.Range("G93").Value = 113.5
$ws.Range("H93").Value = 115.1
$ws.Range("I93").Value = 111.6
$ws.Range("J93").Value = 112.8
$ws.Range("B94").Value = 111.8
$ws.Range("F94").Value = 111
$ws.Range("H94").Value = 116.2
$ws.Range("I94").Value = 111.5
$ws.Range("J94").Value = 113.5
$ws.Range("B95").Value = 113.6
$ws.Range("C95").Value = 105.6
$ws.Range("D95").Value = 97.3
$ws.Range("E95").Value = 104.6
$ws.Range("F95").Value = 112.6
$ws.Range("G95").Value = 116.7
$ws.Range("H95").Value = 117.9
$ws.Range("I95").Value = 113.1
$ws.Range("J95").Value = 115.3
$ws.Range("B96").Value = 114.1
$ws.Range("C96").Value = 107.6
$ws.Range("D96").Value = 102.8
$ws.Range("F96").Value = 113.1
$ws.Range("H96").Value = 118.1
$ws.Range("I96").Value = 113.8
$ws.Range("J96").Value = 115.3
$ws.Range("B97").Value = 109.3
$ws.Range("C97").Value = 105.3
$ws.Range("D97").Value = 102
$ws.Range("E97").Value = 101.7
$ws.Range("F97").Value = 110.5
$ws.Range("G97").Value = 108
$ws.Range("H97").Value = 112.3
$ws.Range("I97").Value = 109.2
$ws.Range("J97").Value = 110
$ws.Range("B98").Value = 111.7
$ws.Range("C98").Value = 106.1
$ws.Range("E98").Value = 102.8
$ws.Range("F98").Value = 112.2
$ws.Range("H98").Value = 115.7
$ws.Range("I98").Value = 111.7
$ws.Range("J98").Value = 112.7
$ws.Range("B99").Value = 98.09999999999999
$ws.Range("C99").Value = 97.2
$ws.Range("D99").Value = 100.1
$ws.Range("E99").Value = 92.59999999999999
$ws.Range("F99").Value = 98.7
$ws.Range("G99").Value = 94.3
$ws.Range("H99").Value = 100.1
$ws.Range("I99").Value = 98.40000000000001
$ws.Range("J99").Value = 97.8
$ws.Range("B100").Value = 102.9
$ws.Range("C100").Value = 98.3
$ws.Range("F100").Value = 94.59999999999999
$ws.Range("H100").Value = 103.3
$ws.Range("I100").Value = 102.4
$ws.Range("J100").Value = 102.9
$ws.Range("B101").Value = 109.1
$ws.Range("C101").Value = 103.6
$ws.Range("D101").Value = 99.40000000000001
$ws.Range("E101").Value = 105.3
$ws.Range("F101").Value = 105.4
$ws.Range("G101").Value = 122.9
$ws.Range("H101").Value = 109.3
$ws.Range("I101").Value = 108.5
$ws.Range("J101").Value = 110.1
$ws.Range("B102").Value = 113.2
$ws.Range("C102").Value = 104.5
$ws.Range("F102").Value = 106.8
$ws.Range("H102").Value = 116.3
$ws.Range("I102").Value = 112.6
$ws.Range("J102").Value = 114.8
$ws.Range("B103").Value = 115.1
$ws.Range("C103").Value = 106.7
$ws.Range("D103").Value = 103
$ws.Range("E103").Value = 107.4
$ws.Range("F103").Value = 109.6
$ws.Range("G103").Value = 135.7
$ws.Range("H103").Value = 115.1
$ws.Range("I103").Value = 113.9
$ws.Range("J103").Value = 116.7
$ws.Range("B104").Value = 121.2
$ws.Range("C104").Value = 107.6
$ws.Range("D104").Value = 98.8
$ws.Range("E104").Value = 111.2
$ws.Range("F104").Value = 113.3
$ws.Range("G104").Value = 145.1
$ws.Range("H104").Value = 123.6
$ws.Range("I104").Value = 119.4
$ws.Range("J104").Value = 124.2
